$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# forNewCustomer sheet: duplicate leather data rows, shift by +5
# ---------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("forNewCustomer")

# Contact numbers C2:C6
$ws3.Cells.Item(2,3).Value = 9881012100
$ws3.Cells.Item(3,3).Value = 9881012101
$ws3.Cells.Item(4,3).Value = 9881012102
$ws3.Cells.Item(5,3).Value = 9881012103
$ws3.Cells.Item(6,3).Value = 9881012104

# Names D2:D6 "abced Test 90..94" -> "abced Test 95..99"
$ws3.Cells.Item(2,4).Value = "abced Test 95"
$ws3.Cells.Item(3,4).Value = "abced Test 96"
$ws3.Cells.Item(4,4).Value = "abced Test 97"
$ws3.Cells.Item(5,4).Value = "abced Test 98"
$ws3.Cells.Item(6,4).Value = "abced Test 99"

# Emails E2:E6 "testuser195..199@mail.com" -> "testuser200..204@mail.com"
$ws3.Cells.Item(2,5).Value = "testuser200@mail.com"
$ws3.Cells.Item(3,5).Value = "testuser201@mail.com"
$ws3.Cells.Item(4,5).Value = "testuser202@mail.com"
$ws3.Cells.Item(5,5).Value = "testuser203@mail.com"
$ws3.Cells.Item(6,5).Value = "testuser204@mail.com"

# Refresh the mail hyperlinks on E3:E6 to point at the new addresses
$links = $ws3.Hyperlinks
$links.Delete()
$ws3.Hyperlinks.Add($ws3.Range("E2"), "mailto:testuser195@mail.com")
$ws3.Hyperlinks.Add($ws3.Range("E3:E6"), "mailto:testuser200@mail.com", [Type]::Missing, [Type]::Missing, "testuser200@mail.com")
$ws3.Hyperlinks.Add($ws3.Range("E3"), "mailto:testuser201@mail.com")
$ws3.Hyperlinks.Add($ws3.Range("E4"), "mailto:testuser202@mail.com")
$ws3.Hyperlinks.Add($ws3.Range("E5"), "mailto:testuser203@mail.com")
$ws3.Hyperlinks.Add($ws3.Range("E6"), "mailto:testuser204@mail.com")

# Update the view: selection moves to D2, no special top-left cell anymore
$ws3.Activate()
$ws3.Range("D2").Select()

# ---------------------------------------------------------------------
# forSync sheet: selection narrows from X2:X6 down to X2, scroll shifts
# ---------------------------------------------------------------------
$ws4 = $wb.Worksheets.Item("forSync")
$ws4.Activate()
$ws4.Range("X2").Select()
$win = $excel.ActiveWindow
$win.ScrollColumn = 3
$win.ScrollRow = 1

# ---------------------------------------------------------------------
# searchInput sheet becomes the active tab
# ---------------------------------------------------------------------
$ws5 = $wb.Worksheets.Item("searchInput")
$ws5.Activate()
